$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7097079753875732
$ws.Range("B1").Value = 2.35954737663269
$ws.Range("C1").Value = 4.994473934173584
$ws.Range("D1").Value = 2.973145723342896
$ws.Range("E1").Value = 0.8116840720176697
